# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the first data
# row (row 2) on both the "zh-cn" and "de-de" worksheets, reflecting a
# fresh report run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-18 20:55:53"
$zhcn.Range("H2").Value = "2016-03-18 20:56:13"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-18 20:55:56"
$dede.Range("H2").Value = "2016-03-18 20:56:22"
